$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20; this pushes the existing rows 20..59
# down to 21..60 (matching the structural shift seen in the diff).
$ws.Rows("20:20").Insert()

# Populate the newly inserted row 20 with the new record's data.
$ws.Range("A20").Value = 4
$ws.Range("B20").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C20").Value = "Los Lagos"
$ws.Range("D20").Value = 44414
$ws.Range("E20").Value = 10
$ws.Range("F20").Value = 100112022
$ws.Range("G20").Value = "Arveja Verde"
$ws.Range("H20").Value = "Perfection"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 70
$ws.Range("K20").Value = 40000
$ws.Range("L20").Value = 40000
$ws.Range("M20").Value = 40000
$ws.Range("N20").Value = "$/malla 25 kilos"
$ws.Range("O20").Value = "Provincia de Huasco"
$ws.Range("P20").Value = 1600
$ws.Range("Q20").Value = 25
$ws.Range("R20").Value = "Hortaliza"
